# Apply the Eterno data updates:
#  - Products: reduce stock for two items (post-sale/order inventory decrement)
#  - POS_Sales: append a new in-store sale record
#  - Customer_Orders: append a new online order record

$wb = $excel.ActiveWorkbook

# --- Products sheet: stock adjustments ---
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Range("E3").Value = 13   # Eterno Grace stock: 15 -> 13
$wsProducts.Range("E10").Value = 9   # Eterno Drift stock: 10 -> 9

# --- POS_Sales sheet: new sale row (row 6) ---
$wsSales = $wb.Worksheets.Item("POS_Sales")
$wsSales.Range("A6").Value = 5
$wsSales.Range("B6").Value = 1
$wsSales.Range("C6").Value = 1498
$wsSales.Range("D6").Value = "cash"
$wsSales.Range("E6").Value = "voucher"
$wsSales.Range("F6").Value = 100
$wsSales.Range("G6").Value = '[{"product_id": 2, "name": "Eterno Grace", "price": 799, "quantity": 2, "stock": 15}]'
$wsSales.Range("H6").Value = "2025-11-09 15:37:44"

# --- Customer_Orders sheet: new order row (row 11) ---
$wsOrders = $wb.Worksheets.Item("Customer_Orders")
$wsOrders.Range("A11").Value = 10
$wsOrders.Range("B11").Value = 3
$wsOrders.Range("C11").Value = "kaizen"
$wsOrders.Range("D11").Value = "boarratjabol@gmail.com"
$wsOrders.Range("E11").Value = "blk 32 lot 12, Paloma Street, Garden of Edem, Greece, Europe"
$wsOrders.Range("F11").Value = 3999
$wsOrders.Range("G11").Value = 83
$wsOrders.Range("H11").Value = 3982
$wsOrders.Range("I11").Value = "cod"
$wsOrders.Range("J11").Value = "completed"
$wsOrders.Range("K11").Value = '[{"product_id": 9, "product_name": "Eterno Drift", "quantity": 1, "price": 3999.0}]'
$wsOrders.Range("L11").Value = "2025-11-09 15:36:22"
